$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color constants used by the PASS/FAIL highlighting in this report
# (Excel Long color values, computed from the RGB hex used in the workbook)
$greenFont  = 32768     # FF008000
$greenFill  = 9498256   # FF90EE90
$redFont    = 255       # FFFF0000
$redFill    = 12695295  # FFFFB6C1

# ---------------------------------------------------------------------------
# KDT01 (row 2): Actual result text now differs from expected -> test FAILs
# ---------------------------------------------------------------------------
$ws.Range("H2").Value = '{"statusCode":500,"response":{"errorMessage":"Lịch hẹn không tồn tại!"}}'
$ws.Range("J2").Value = "FAIL"
$ws.Range("K2").Value = "23.78ms"

# ---------------------------------------------------------------------------
# KDT02 (row 3): still PASS, only timing changes
# ---------------------------------------------------------------------------
$ws.Range("K3").Value = "2.58ms"

# ---------------------------------------------------------------------------
# KDT03 (row 4): Expected result message reworded -> now mismatches actual -> FAIL
# ---------------------------------------------------------------------------
$ws.Range("H4").Value = '{"statusCode":403,"response":{"errorMessage":"Quyền truy cập vào lịch hẹn bị từ chối!"}}'
$ws.Range("J4").Value = "FAIL"
$ws.Range("K4").Value = "1.01ms"

# ---------------------------------------------------------------------------
# KDT04 (row 5): Actual result no longer includes the "note" field -> now matches expected -> PASS, timing changes
# ---------------------------------------------------------------------------
$ws.Range("I5").Value = '{"statusCode":200,"response":{"message":"Tạo toa thuốc thành công!"}}'
$ws.Range("K5").Value = "18.50ms"

# ---------------------------------------------------------------------------
# KDT05 (row 6): still PASS, only timing changes
# ---------------------------------------------------------------------------
$ws.Range("K6").Value = "2.78ms"

# ---------------------------------------------------------------------------
# KDT06 (row 7): still FAIL, only timing changes
# ---------------------------------------------------------------------------
$ws.Range("K7").Value = "0.82ms"

# ---------------------------------------------------------------------------
# KDT07 (row 8): still PASS, only timing changes
# ---------------------------------------------------------------------------
$ws.Range("K8").Value = "1.81ms"

# ---------------------------------------------------------------------------
# KDT08 (row 9): still FAIL, only timing changes
# ---------------------------------------------------------------------------
$ws.Range("K9").Value = "1.14ms"

# ---------------------------------------------------------------------------
# KDT09 (row 10): still FAIL, only timing changes
# ---------------------------------------------------------------------------
$ws.Range("K10").Value = "0.62ms"

# ---------------------------------------------------------------------------
# KDT10 (row 11): still PASS, only timing changes
# ---------------------------------------------------------------------------
$ws.Range("K11").Value = "1.15ms"

# ---------------------------------------------------------------------------
# Re-apply PASS (green) / FAIL (red) colouring for every status cell so the
# colour always matches the (possibly updated) text.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Range("J$r")
    $status = $cell.Value2
    if ($status -eq "PASS") {
        $cell.Font.Color = $greenFont
        $cell.Interior.Color = $greenFill
    } else {
        $cell.Font.Color = $redFont
        $cell.Interior.Color = $redFill
    }
}
